$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 64, shifting existing rows 64-69 down to 65-70
$ws.Rows.Item(64).Insert()

# Populate the new row 64 with the new record
$ws.Cells.Item(64, 1).Value = 3
$ws.Cells.Item(64, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(64, 3).Value = "Coquimbo"
$ws.Cells.Item(64, 4).Value = 44783
$ws.Cells.Item(64, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(64, 5).Value = 5
$ws.Cells.Item(64, 6).Value = 100112035
$ws.Cells.Item(64, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 50
$ws.Cells.Item(64, 11).Value = 15000
$ws.Cells.Item(64, 12).Value = 15000
$ws.Cells.Item(64, 13).Value = 15000
$ws.Cells.Item(64, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(64, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(64, 16).Value = 1000
$ws.Cells.Item(64, 17).Value = 15
$ws.Cells.Item(64, 18).Value = "Hortaliza"
